$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 46073
$ws.Range("C3").Value = 46073
$ws.Range("C4").Value = 46073
$ws.Range("C5").Value = 46073
$ws.Range("C6").Value = 46073
$ws.Range("C7").Value = 46073
$ws.Range("C8").Value = 46073
$ws.Range("C9").Value = 46073
$ws.Range("C10").Value = 46073
$ws.Range("A11").Value = 'A 27636-2023'
$ws.Range("B11").Value = 45097
$ws.Range("C11").Value = 46073
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/artfynd/A 27636-2023 artfynd.xlsx", "A 27636-2023")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/kartor/A 27636-2023 karta.png", "A 27636-2023")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomål/A 27636-2023 FSC-klagomål.docx", "A 27636-2023")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomålsmail/A 27636-2023 FSC-klagomål mail.docx", "A 27636-2023")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsyn/A 27636-2023 tillsynsbegäran.docx", "A 27636-2023")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsynsmail/A 27636-2023 tillsynsbegäran mail.docx", "A 27636-2023")'
$ws.Range("A12").Value = 'A 38039-2022'
$ws.Range("B12").Value = 44811
$ws.Range("C12").Value = 46073
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/artfynd/A 38039-2022 artfynd.xlsx", "A 38039-2022")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/kartor/A 38039-2022 karta.png", "A 38039-2022")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomål/A 38039-2022 FSC-klagomål.docx", "A 38039-2022")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomålsmail/A 38039-2022 FSC-klagomål mail.docx", "A 38039-2022")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsyn/A 38039-2022 tillsynsbegäran.docx", "A 38039-2022")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsynsmail/A 38039-2022 tillsynsbegäran mail.docx", "A 38039-2022")'
$ws.Range("C13").Value = 46073
$ws.Range("C14").Value = 46073
$ws.Range("C15").Value = 46073
$ws.Range("C16").Value = 46073
$ws.Range("C17").Value = 46073
$ws.Range("A18").Value = 'A 2864-2026'
$ws.Range("B18").Value = 46038
$ws.Range("C18").Value = 46073
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 1
$ws.Range("L18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = 'Spillkråka' + "`r`n" + 'Mindre märgborre' + "`r`n" + 'Blåsippa'
$ws.Range("S18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/artfynd/A 2864-2026 artfynd.xlsx", "A 2864-2026")'
$ws.Range("T18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/kartor/A 2864-2026 karta.png", "A 2864-2026")'
$ws.Range("V18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomål/A 2864-2026 FSC-klagomål.docx", "A 2864-2026")'
$ws.Range("W18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomålsmail/A 2864-2026 FSC-klagomål mail.docx", "A 2864-2026")'
$ws.Range("X18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsyn/A 2864-2026 tillsynsbegäran.docx", "A 2864-2026")'
$ws.Range("Y18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsynsmail/A 2864-2026 tillsynsbegäran mail.docx", "A 2864-2026")'
$ws.Range("Z18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/fåglar/A 2864-2026 prioriterade fågelarter.docx", "A 2864-2026")'
$ws.Range("A19").Value = 'A 61963-2025'
$ws.Range("B19").Value = 46003
$ws.Range("C19").Value = 46073
$ws.Range("G19").Value = 1.1
$ws.Range("S19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/artfynd/A 61963-2025 artfynd.xlsx", "A 61963-2025")'
$ws.Range("T19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/kartor/A 61963-2025 karta.png", "A 61963-2025")'
$ws.Range("V19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomål/A 61963-2025 FSC-klagomål.docx", "A 61963-2025")'
$ws.Range("W19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomålsmail/A 61963-2025 FSC-klagomål mail.docx", "A 61963-2025")'
$ws.Range("X19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsyn/A 61963-2025 tillsynsbegäran.docx", "A 61963-2025")'
$ws.Range("Y19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsynsmail/A 61963-2025 tillsynsbegäran mail.docx", "A 61963-2025")'
$ws.Range("Z19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/fåglar/A 61963-2025 prioriterade fågelarter.docx", "A 61963-2025")'
$ws.Range("A20").Value = 'A 35242-2024'
$ws.Range("B20").Value = 45530.55440972222
$ws.Range("C20").Value = 46073
$ws.Range("G20").Value = 0.9
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 1
$ws.Range("P20").Value = 1
$ws.Range("R20").Value = 'Ryl' + "`r`n" + 'Grönpyrola' + "`r`n" + 'Skogsknipprot'
$ws.Range("S20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/artfynd/A 35242-2024 artfynd.xlsx", "A 35242-2024")'
$ws.Range("T20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/kartor/A 35242-2024 karta.png", "A 35242-2024")'
$ws.Range("V20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomål/A 35242-2024 FSC-klagomål.docx", "A 35242-2024")'
$ws.Range("W20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomålsmail/A 35242-2024 FSC-klagomål mail.docx", "A 35242-2024")'
$ws.Range("X20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsyn/A 35242-2024 tillsynsbegäran.docx", "A 35242-2024")'
$ws.Range("Y20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsynsmail/A 35242-2024 tillsynsbegäran mail.docx", "A 35242-2024")'
$ws.Range("Z20").ClearContents()
$ws.Range("A21").Value = 'A 15600-2025'
$ws.Range("B21").Value = 45747
$ws.Range("C21").Value = 46073
$ws.Range("G21").Value = 1.1
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 2
$ws.Range("R21").Value = 'Murgröna' + "`r`n" + 'Scharlakansvårskål agg.'
$ws.Range("S21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/artfynd/A 15600-2025 artfynd.xlsx", "A 15600-2025")'
$ws.Range("T21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/kartor/A 15600-2025 karta.png", "A 15600-2025")'
$ws.Range("V21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomål/A 15600-2025 FSC-klagomål.docx", "A 15600-2025")'
$ws.Range("W21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomålsmail/A 15600-2025 FSC-klagomål mail.docx", "A 15600-2025")'
$ws.Range("X21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsyn/A 15600-2025 tillsynsbegäran.docx", "A 15600-2025")'
$ws.Range("Y21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsynsmail/A 15600-2025 tillsynsbegäran mail.docx", "A 15600-2025")'
$ws.Range("A22").Value = 'A 62231-2023'
$ws.Range("B22").Value = 45267
$ws.Range("C22").Value = 46073
$ws.Range("G22").Value = 1.6
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("R22").Value = 'Murgröna' + "`r`n" + 'Blåsippa'
$ws.Range("S22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/artfynd/A 62231-2023 artfynd.xlsx", "A 62231-2023")'
$ws.Range("T22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/kartor/A 62231-2023 karta.png", "A 62231-2023")'
$ws.Range("V22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomål/A 62231-2023 FSC-klagomål.docx", "A 62231-2023")'
$ws.Range("W22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomålsmail/A 62231-2023 FSC-klagomål mail.docx", "A 62231-2023")'
$ws.Range("X22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsyn/A 62231-2023 tillsynsbegäran.docx", "A 62231-2023")'
$ws.Range("Y22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsynsmail/A 62231-2023 tillsynsbegäran mail.docx", "A 62231-2023")'
$ws.Range("A23").Value = 'A 40361-2023'
$ws.Range("B23").Value = 45169
$ws.Range("C23").Value = 46073
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1
$ws.Range("M23").Value = 1
$ws.Range("O23").Value = 2
$ws.Range("P23").Value = 1
$ws.Range("R23").Value = 'Lundalm' + "`r`n" + 'Ängsskära'
$ws.Range("S23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/artfynd/A 40361-2023 artfynd.xlsx", "A 40361-2023")'
$ws.Range("T23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/kartor/A 40361-2023 karta.png", "A 40361-2023")'
$ws.Range("V23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomål/A 40361-2023 FSC-klagomål.docx", "A 40361-2023")'
$ws.Range("W23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/klagomålsmail/A 40361-2023 FSC-klagomål mail.docx", "A 40361-2023")'
$ws.Range("X23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsyn/A 40361-2023 tillsynsbegäran.docx", "A 40361-2023")'
$ws.Range("Y23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0840/tillsynsmail/A 40361-2023 tillsynsbegäran mail.docx", "A 40361-2023")'
$ws.Range("C24").Value = 46073
$ws.Range("C25").Value = 46073
$ws.Range("C26").Value = 46073
$ws.Range("C27").Value = 46073
$ws.Range("C28").Value = 46073
$ws.Range("C29").Value = 46073
$ws.Range("C30").Value = 46073
$ws.Range("A31").Value = 'A 13040-2024'
$ws.Range("B31").Value = 45385
$ws.Range("C31").Value = 46073
$ws.Range("G31").Value = 1.2
$ws.Range("A32").Value = 'A 13384-2023'
$ws.Range("B32").Value = 45005
$ws.Range("C32").Value = 46073
$ws.Range("G32").Value = 0.9
$ws.Range("A33").Value = 'A 57407-2025'
$ws.Range("B33").Value = 45980.44351851852
$ws.Range("C33").Value = 46073
$ws.Range("A34").Value = 'A 62183-2025'
$ws.Range("B34").Value = 46006.5356712963
$ws.Range("C34").Value = 46073
$ws.Range("G34").Value = 0.8
$ws.Range("A35").Value = 'A 5733-2026'
$ws.Range("B35").Value = 46051.5745949074
$ws.Range("C35").Value = 46073
$ws.Range("G35").Value = 1.1
$ws.Range("A36").Value = 'A 63676-2023'
$ws.Range("B36").Value = 45275.64739583333
$ws.Range("C36").Value = 46073
$ws.Range("G36").Value = 1
$ws.Range("C37").Value = 46073
